$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "DOI"
$ws.Range("B1").Value = "Citations"
$ws.Range("A2").Value = 10.276199999999999
$ws.Range("B2").Value = 66
$ws.Range("A3").Value = "10.33322l"
$ws.Range("B3").Value = 1

$ws.Range("E5").Select() | Out-Null
